$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2161.9
$ws.Range("J17").Value = 2291
$ws.Range("L17").Value = 6873
$ws.Range("N17").Value = -7209

$ws.Range("H40").Value = 2587.75
$ws.Range("I40").Value = 1783.6666
$ws.Range("K40").Value = 1783.6666
$ws.Range("M40").Value = -1608.6666

$ws.Range("H43").Value = 7604.0586
$ws.Range("I43").Value = 6571.2856
$ws.Range("K43").Value = 6571.2856
$ws.Range("M43").Value = -6502.2856

$ws.Range("H58").Value = 316.66666
$ws.Range("I58").Value = 316.66666
$ws.Range("K58").Value = 949.9999799999999
$ws.Range("M58").Value = -799.9999799999999

$ws.Range("H76").Value = 3485
$ws.Range("I76").Value = 3147.9167
$ws.Range("K76").Value = 3147.9167
$ws.Range("M76").Value = -2832.9167

$ws.Range("H79").Value = 3485
$ws.Range("I79").Value = 3147.9167
$ws.Range("K79").Value = 3147.9167
$ws.Range("M79").Value = -2055.9167

$ws.Range("H80").Value = 723.8077
$ws.Range("J80").Value = 814.0625
$ws.Range("L80").Value = 2442.1875
$ws.Range("N80").Value = -4438.1875

$ws.Range("H83").Value = 723.8077
$ws.Range("J83").Value = 814.0625
$ws.Range("L83").Value = 7326.5625
$ws.Range("N83").Value = -17310.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1877.5781
$ws.Range("I61").Value = 1877.5781
$ws.Range("K61").Value = 1877.5781
$ws.Range("M61").Value = -1665.5781

$ws.Range("H110").Value = 584.8333
$ws.Range("I110").Value = 461.8
$ws.Range("K110").Value = 461.8
$ws.Range("M110").Value = 1583.2

$ws.Range("H122").Value = 2093
$ws.Range("I122").Value = 1657.3334
$ws.Range("K122").Value = 4972.0002
$ws.Range("M122").Value = -2522.0002

$ws.Range("H136").Value = 1877.5781
$ws.Range("I136").Value = 1877.5781
$ws.Range("K136").Value = 5632.7343
$ws.Range("M136").Value = -3082.7343

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 35000
$ws.Range("I9").Value = 35000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 35000
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("M9").Value = -34832

$ws.Range("H107").Value = 2932.8064
$ws.Range("I107").Value = 1469.421
$ws.Range("J107").Value = 5249.8335
$ws.Range("K107").Value = 1469.421
$ws.Range("L107").Value = 5249.8335
$ws.Range("M107").Value = 450.579
$ws.Range("N107").Value = -9089.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5625.75
$ws.Range("J62").Value = 6928.5713
$ws.Range("L62").Value = 6928.5713
$ws.Range("N62").Value = -8176.5713

$ws.Range("H65").Value = 5625.75
$ws.Range("J65").Value = 6928.5713
$ws.Range("L65").Value = 34642.85649999999
$ws.Range("N65").Value = -40882.85649999999

$ws.Range("H107").Value = 365.57144
$ws.Range("I107").Value = 385
$ws.Range("K107").Value = 385
$ws.Range("M107").Value = 1535

$ws.Range("H141").Value = 290163
$ws.Range("J141").Value = 290163
$ws.Range("L141").Value = 290163
$ws.Range("N141").Value = -300523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1288823.2
$ws.Range("I4").Value = 27697.21
$ws.Range("J4").Value = 25250218
$ws.Range("K4").Value = 83091.63
$ws.Range("L4").Value = 75750654
$ws.Range("M4").Value = -82979.63
$ws.Range("N4").Value = -75750878

$ws.Range("H32").Value = 2634.5
$ws.Range("I32").Value = 1915
$ws.Range("K32").Value = 5745
$ws.Range("M32").Value = -5462

$ws.Range("H62").Value = 10500
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372

$ws.Range("H65").Value = 10500
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864

$ws.Range("H120").Value = 66628.3
$ws.Range("I120").Value = 101418.25
$ws.Range("J120").Value = 43435
$ws.Range("K120").Value = 304254.75
$ws.Range("L120").Value = 130305
$ws.Range("M120").Value = -299416.75
$ws.Range("N120").Value = -139981

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3496.5
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 3496.5
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984

$ws.Range("H113").Value = 5911.077
$ws.Range("I113").Value = 3495.2856
$ws.Range("K113").Value = 3495.2856
$ws.Range("M113").Value = -1325.2856

$ws.Range("H126").Value = 2520.389
$ws.Range("I126").Value = 1995.3
$ws.Range("K126").Value = 5985.9
$ws.Range("M126").Value = -3515.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4749
$ws.Range("I7").Value = 4749
$ws.Range("K7").Value = 4749
$ws.Range("M7").Value = -4637

$ws.Range("H46").Value = 2564.7693
$ws.Range("I46").Value = 1019.8
$ws.Range("K46").Value = 1019.8
$ws.Range("M46").Value = -831.8

$ws.Range("H122").Value = 5275.727
$ws.Range("I122").Value = 5750.6113
$ws.Range("K122").Value = 17251.8339
$ws.Range("M122").Value = -14801.8339

$ws.Range("H126").Value = 4749
$ws.Range("I126").Value = 4749
$ws.Range("K126").Value = 14247
$ws.Range("M126").Value = -11777

$ws.Range("H132").Value = 3720.08
$ws.Range("I132").Value = 3571.5715
$ws.Range("K132").Value = 10714.7145
$ws.Range("M132").Value = -8184.7145

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2749.6843
$ws.Range("I81").Value = 1885.3077
$ws.Range("K81").Value = 3770.6154
$ws.Range("M81").Value = -2709.6154

$ws.Range("H84").Value = 2749.6843
$ws.Range("I84").Value = 1885.3077
$ws.Range("K84").Value = 18853.077
$ws.Range("M84").Value = -13549.077

$ws.Range("H119").Value = 66666
$ws.Range("J119").Value = 66666
$ws.Range("L119").Value = 66666
$ws.Range("N119").Value = -76342

$ws.Range("H126").Value = 1593.5714
$ws.Range("I126").Value = 1620
$ws.Range("J126").Value = 1527.5
$ws.Range("K126").Value = 4860
$ws.Range("L126").Value = 4582.5
$ws.Range("M126").Value = -2390
$ws.Range("N126").Value = -9522.5
